$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Sheet "总计": update existing quarter row to Q4 values, then
#    append a new row restoring the old Q3 values underneath it.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Copy the formatting of row 2's first cell onto the new row 3 cell
# before writing into it, so it picks up the same style (bold/border).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 0.01

# ------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q3" sheet to create the new
#    "2022-Q4" sheet right in front of it, then overwrite its
#    figures with the Q4 numbers.
# ------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3, $null)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

$wsQ4.Range("D2").Value = "'0.21"
$wsQ4.Range("E2").Value = "'86.58"
$wsQ4.Range("F2").Value = "'3.54"
$wsQ4.Range("G2").Value = "'0.0074"
$wsQ4.Range("H2").Value = 10

$wsQ4.Range("D3").Value = "'0.16"
$wsQ4.Range("E3").Value = "'86.58"
$wsQ4.Range("F3").Value = "'3.54"
$wsQ4.Range("G3").Value = "'0.0057"
$wsQ4.Range("H3").Value = 10

# Writing a leading-apostrophe string stores it as text but also flips on
# the "stored as text" quote-prefix marker on the cell style; strip that
# back off so the cells end up with the same plain/default styling as the
# rest of the numeric-looking text column (matches the source sheet).
$wsQ4.Range("D2:G3").ClearFormats()

Write-Output ("Worksheets: " + $wb.Worksheets.Count)
